$d = $word.ActiveDocument

# Move to the very end of the document content (after "...both locked and lockless",
# before the _GoBack bookmark) and insert two new paragraphs.
$end = $d.Content
$end.Collapse(0)   # wdCollapseEnd
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.MoveStart(1, 1)  # wdCharacter

$end.InsertAfter("26/02/14")
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.MoveStart(1, 1)

$end.InsertAfter("Gathered data from linked list using perf")

$d.Save()
